# Update ST_PROD_01 (row 2)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "1. Điều hướng đến Quản lý SP -> Thêm mới`n2. Điền form thông tin SP`n3. Lưu & Tìm ở trang cuối"
$ws.Range("D2").Value = "Name: Ao Test Auto 1764923200090, Price: 500000"
$ws.Range("E2").Value = "Sản phẩm mới xuất hiện trong danh sách (có hỗ trợ phân trang)"

# Update ST_PROD_02 (row 3)
$ws.Range("C3").Value = "1. Điều hướng đến trang DS`n2. Nhập tên SP vừa tạo vào ô tìm kiếm`n3. Check kết quả"
$ws.Range("D3").Value = "Keyword: Ao Test Auto 1764923200090"
$ws.Range("E3").Value = "Hiển thị đúng sản phẩm (chỉ 1 kết quả)"

# Update ST_PROD_03 (row 4)
$ws.Range("C4").Value = "1. Tìm SP`n2. Click Sửa`n3. Đổi giá thành 1000000 -> Lưu"
$ws.Range("E4").Value = "Giá tiền cập nhật thành công (1,000,000)"

# Update ST_PROD_04 (row 5)
$ws.Range("D5").Value = "Target: Ao Test Auto 1764923200090"
$ws.Range("E5").Value = "Sản phẩm biến mất khỏi danh sách"

# Widen columns C, D, E to fit the new, longer text (mirrors Excel's bestFit autosize
# recalculation that happens after the cell content grows)
$ws.Columns.Item(3).ColumnWidth = 38.833333333333336
$ws.Columns.Item(4).ColumnWidth = 46.333333333333336
$ws.Columns.Item(5).ColumnWidth = 57.0

# Re-fit the row heights (setting multi-line values can mark rows with a custom
# height); AutoFit puts them back to the sheet's default row height, same as the source
$ws.Rows.Item(2).AutoFit() | Out-Null
$ws.Rows.Item(3).AutoFit() | Out-Null
$ws.Rows.Item(4).AutoFit() | Out-Null
